$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old UUID-like values with the new ones, per the diff.
$ws.Range("A2").Value = "8091a8cf-0e6c-4381-b2ea-0e17a9e90b4b"
$ws.Range("A3").Value = "46ef86ea-22cb-48e8-a0b4-a57bd08ccf46"
$ws.Range("G4").Value = "3f391fe9-224c-4a58-a21e-030e759126f5"
$ws.Range("G5").Value = "06f8bec4-4ee1-4926-92db-b5b23e1cb843"
$ws.Range("A6").Value = "2d977ddb-7cd4-4d81-af99-ae58e7a6e3d3"
$ws.Range("G6").Value = "2d977ddb-7cd4-4d81-af99-ae58e7a6e3d3"
$ws.Range("G7").Value = "6997738c-e7db-4bd9-9623-552bde6b9784"
